$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 33.43904566666667
$ws.Range("H2").Value = 100.317137
$ws.Range("I2").Value = 0.5629652610385096
$ws.Range("J2").Value = 0.5629652610385097
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1451143333333333
$ws.Range("N2").Value = 0.435343
$ws.Range("O2").Value = 0.140827770705994
$ws.Range("P2").Value = 0.1408277707059941
$ws.Range("Q2").Value = 4.852484819221221
$ws.Range("R2").Value = 43.672363372991
$ws.Range("S2").Value = 0.07928114269697131
$ws.Range("T2").Value = 0.07928114269697134

$ws.Range("G3").Value = 33.43904566666667
$ws.Range("H3").Value = 100.317137
$ws.Range("I3").Value = 0.5629652610385096
$ws.Range("J3").Value = 0.5629652610385097
$ws.Range("O3").Value = 0.0422072807203407
$ws.Range("P3").Value = 0.0422072807203407
$ws.Range("Q3").Value = 1.454330974134667
$ws.Range("R3").Value = 13.088978767212
$ws.Range("S3").Value = 0.02376123280845226
$ws.Range("T3").Value = 0.02376123280845226

$ws.Range("G4").Value = 33.43904566666667
$ws.Range("H4").Value = 100.317137
$ws.Range("I4").Value = 0.5629652610385096
$ws.Range("J4").Value = 0.5629652610385097
$ws.Range("M4").Value = 0.841832
$ws.Range("N4").Value = 2.525496
$ws.Range("O4").Value = 0.8169649485736653
$ws.Range("P4").Value = 0.8169649485736653
$ws.Range("Q4").Value = 28.15005869166133
$ws.Range("R4").Value = 253.350528224952
$ws.Range("S4").Value = 0.4599228855330861
$ws.Range("T4").Value = 0.4599228855330862

$ws.Range("I5").Value = 0.3606447218168795
$ws.Range("J5").Value = 0.3606447218168795
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1451143333333333
$ws.Range("N5").Value = 0.435343
$ws.Range("O5").Value = 0.140827770705994
$ws.Range("P5").Value = 0.1408277707059941
$ws.Range("Q5").Value = 3.108580864333222
$ws.Range("R5").Value = 27.977227778999
$ws.Range("S5").Value = 0.05078879219035451
$ws.Range("T5").Value = 0.05078879219035452

$ws.Range("I6").Value = 0.3606447218168795
$ws.Range("J6").Value = 0.3606447218168795
$ws.Range("O6").Value = 0.0422072807203407
$ws.Range("P6").Value = 0.0422072807203407
$ws.Range("S6").Value = 0.01522183301403421
$ws.Range("T6").Value = 0.01522183301403422

$ws.Range("I7").Value = 0.3606447218168795
$ws.Range("J7").Value = 0.3606447218168795
$ws.Range("M7").Value = 0.841832
$ws.Range("N7").Value = 2.525496
$ws.Range("O7").Value = 0.8169649485736653
$ws.Range("P7").Value = 0.8169649485736653
$ws.Range("Q7").Value = 18.03338640692533
$ws.Range("R7").Value = 162.300477662328
$ws.Range("S7").Value = 0.2946340966124908
$ws.Range("T7").Value = 0.2946340966124908

$ws.Range("G8").Value = 0.5400056666666667
$ws.Range("H8").Value = 1.620017
$ws.Range("I8").Value = 0.009091301053496209
$ws.Range("J8").Value = 0.009091301053496209
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1451143333333333
$ws.Range("N8").Value = 0.435343
$ws.Range("O8").Value = 0.140827770705994
$ws.Range("P8").Value = 0.1408277707059941
$ws.Range("Q8").Value = 0.07836256231455556
$ws.Range("R8").Value = 0.705263060831
$ws.Range("S8").Value = 0.001280307660180926
$ws.Range("T8").Value = 0.001280307660180926

$ws.Range("G9").Value = 0.5400056666666667
$ws.Range("H9").Value = 1.620017
$ws.Range("I9").Value = 0.009091301053496209
$ws.Range("J9").Value = 0.009091301053496209
$ws.Range("O9").Value = 0.0422072807203407
$ws.Range("P9").Value = 0.0422072807203407
$ws.Range("Q9").Value = 0.02348592645466667
$ws.Range("R9").Value = 0.211373338092
$ws.Range("S9").Value = 0.0003837190956780436
$ws.Range("T9").Value = 0.0003837190956780436

$ws.Range("G10").Value = 0.5400056666666667
$ws.Range("H10").Value = 1.620017
$ws.Range("I10").Value = 0.009091301053496209
$ws.Range("J10").Value = 0.009091301053496209
$ws.Range("M10").Value = 0.841832
$ws.Range("N10").Value = 2.525496
$ws.Range("O10").Value = 0.8169649485736653
$ws.Range("P10").Value = 0.8169649485736653
$ws.Range("Q10").Value = 0.4545940503813334
$ws.Range("R10").Value = 4.091346453432
$ws.Range("S10").Value = 0.00742727429763724
$ws.Range("T10").Value = 0.00742727429763724

$ws.Range("G11").Value = 2.881522666666667
$ws.Range("H11").Value = 8.644568
$ws.Range("I11").Value = 0.04851206509895859
$ws.Range("J11").Value = 0.04851206509895859
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.1451143333333333
$ws.Range("N11").Value = 0.435343
$ws.Range("O11").Value = 0.140827770705994
$ws.Range("P11").Value = 0.1408277707059941
$ws.Range("Q11").Value = 0.4181502407582222
$ws.Range("R11").Value = 3.763352166824
$ws.Range("S11").Value = 0.006831845980230396
$ws.Range("T11").Value = 0.006831845980230397

$ws.Range("G12").Value = 2.881522666666667
$ws.Range("H12").Value = 8.644568
$ws.Range("I12").Value = 0.04851206509895859
$ws.Range("J12").Value = 0.04851206509895859
$ws.Range("O12").Value = 0.0422072807203407
$ws.Range("P12").Value = 0.0422072807203407
$ws.Range("Q12").Value = 0.1253231838186667
$ws.Range("R12").Value = 1.127908654368
$ws.Range("S12").Value = 0.002047562349955188
$ws.Range("T12").Value = 0.002047562349955188

$ws.Range("G13").Value = 2.881522666666667
$ws.Range("H13").Value = 8.644568
$ws.Range("I13").Value = 0.04851206509895859
$ws.Range("J13").Value = 0.04851206509895859
$ws.Range("M13").Value = 0.841832
$ws.Range("N13").Value = 2.525496
$ws.Range("O13").Value = 0.8169649485736653
$ws.Range("P13").Value = 0.8169649485736653
$ws.Range("Q13").Value = 2.425757989525334
$ws.Range("R13").Value = 21.831821905728
$ws.Range("S13").Value = 0.03963265676877301
$ws.Range("T13").Value = 0.03963265676877301

$ws.Range("G14").Value = 0.4104343333333333
$ws.Range("H14").Value = 1.231303
$ws.Range("I14").Value = 0.006909894316586209
$ws.Range("J14").Value = 0.006909894316586209
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1451143333333333
$ws.Range("N14").Value = 0.435343
$ws.Range("O14").Value = 0.140827770705994
$ws.Range("P14").Value = 0.1408277707059941
$ws.Range("Q14").Value = 0.05955990465877777
$ws.Range("R14").Value = 0.536039141929
$ws.Range("S14").Value = 0.000973105012418854
$ws.Range("T14").Value = 0.0009731050124188542

$ws.Range("G15").Value = 0.4104343333333333
$ws.Range("H15").Value = 1.231303
$ws.Range("I15").Value = 0.006909894316586209
$ws.Range("J15").Value = 0.006909894316586209
$ws.Range("O15").Value = 0.0422072807203407
$ws.Range("P15").Value = 0.0422072807203407
$ws.Range("Q15").Value = 0.01785061002533334
$ws.Range("R15").Value = 0.160655490228
$ws.Range("S15").Value = 0.0002916478491680409
$ws.Range("T15").Value = 0.0002916478491680409

$ws.Range("G16").Value = 0.4104343333333333
$ws.Range("H16").Value = 1.231303
$ws.Range("I16").Value = 0.006909894316586209
$ws.Range("J16").Value = 0.006909894316586209
$ws.Range("M16").Value = 0.841832
$ws.Range("N16").Value = 2.525496
$ws.Range("O16").Value = 0.8169649485736653
$ws.Range("P16").Value = 0.8169649485736653
$ws.Range("Q16").Value = 0.3455167556986667
$ws.Range("R16").Value = 3.109650801288
$ws.Range("S16").Value = 0.005645141454999314
$ws.Range("T16").Value = 0.005645141454999314

$ws.Range("G17").Value = 0.7054563333333334
$ws.Range("H17").Value = 2.116369
$ws.Range("I17").Value = 0.0118767566755699
$ws.Range("J17").Value = 0.0118767566755699
$ws.Range("K17").Value = 1
$ws.Range("L17").Value = 0.3333333333333333
$ws.Range("M17").Value = 0.1451143333333333
$ws.Range("N17").Value = 0.435343
$ws.Range("O17").Value = 0.140827770705994
$ws.Range("P17").Value = 0.1408277707059941
$ws.Range("Q17").Value = 0.1023718255074444
$ws.Range("R17").Value = 0.9213464295670001
$ws.Range("S17").Value = 0.001672577165838041
$ws.Range("T17").Value = 0.001672577165838042

$ws.Range("G18").Value = 0.7054563333333334
$ws.Range("H18").Value = 2.116369
$ws.Range("I18").Value = 0.0118767566755699
$ws.Range("J18").Value = 0.0118767566755699
$ws.Range("O18").Value = 0.0422072807203407
$ws.Range("P18").Value = 0.0422072807203407
$ws.Range("Q18").Value = 0.03068170684933334
$ws.Range("R18").Value = 0.276135361644
$ws.Range("S18").Value = 0.0005012856030529589
$ws.Range("T18").Value = 0.0005012856030529591

$ws.Range("G19").Value = 0.7054563333333334
$ws.Range("H19").Value = 2.116369
$ws.Range("I19").Value = 0.0118767566755699
$ws.Range("J19").Value = 0.0118767566755699
$ws.Range("M19").Value = 0.841832
$ws.Range("N19").Value = 2.525496
$ws.Range("O19").Value = 0.8169649485736653
$ws.Range("P19").Value = 0.8169649485736653
$ws.Range("Q19").Value = 0.5938757160026668
$ws.Range("R19").Value = 5.344881444024001
$ws.Range("S19").Value = 0.009702893906678896
$ws.Range("T19").Value = 0.009702893906678898
